$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 10165.667
$ws.Range("I51").Value = 9874
$ws.Range("J51").Value = 10399
$ws.Range("K51").Value = 9874
$ws.Range("L51").Value = 10399
$ws.Range("M51").Value = -9390
$ws.Range("N51").Value = -11367
$ws.Range("H64").Value = 5410
$ws.Range("I64").Value = 5732
$ws.Range("K64").Value = 5732
$ws.Range("M64").Value = -5484
$ws.Range("H67").Value = 5410
$ws.Range("I67").Value = 5732
$ws.Range("K67").Value = 5732
$ws.Range("M67").Value = -4874
$ws.Range("H70").Value = 1301.8
$ws.Range("I70").Value = 1291.2858
$ws.Range("J70").Value = 1311
$ws.Range("K70").Value = 3873.8574
$ws.Range("L70").Value = 3933
$ws.Range("M70").Value = -3603.8574
$ws.Range("N70").Value = -4473
$ws.Range("H73").Value = 1301.8
$ws.Range("I73").Value = 1291.2858
$ws.Range("J73").Value = 1311
$ws.Range("K73").Value = 3873.8574
$ws.Range("L73").Value = 3933
$ws.Range("M73").Value = -2937.8574
$ws.Range("N73").Value = -5805
$ws.Range("H98").Value = 729.2857
$ws.Range("I98").Value = 729.2857
$ws.Range("K98").Value = 729.2857
$ws.Range("M98").Value = 768.7143
$ws.Range("H100").Value = 2637.2
$ws.Range("J100").Value = 2452.5
$ws.Range("L100").Value = 2452.5
$ws.Range("N100").Value = -3534.5
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").ClearContents()
$ws.Range("N105").Value = 0
$ws.Range("H122").Value = 729.2857
$ws.Range("I122").Value = 729.2857
$ws.Range("K122").Value = 2187.8571
$ws.Range("M122").Value = 262.1428999999998

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12553.8
$ws.Range("I32").Value = 9260.727999999999
$ws.Range("J32").Value = 36703
$ws.Range("K32").Value = 9260.727999999999
$ws.Range("L32").Value = 36703
$ws.Range("M32").Value = -8973.727999999999
$ws.Range("N32").Value = -37277
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").ClearContents()
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = 0

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 4181.6
$ws.Range("I19").Value = 225.5
$ws.Range("J19").Value = 20006
$ws.Range("K19").Value = 225.5
$ws.Range("L19").Value = 20006
$ws.Range("M19").Value = -55.5
$ws.Range("N19").Value = -20346
$ws.Range("H22").Value = 386.66666
$ws.Range("I22").Value = 308
$ws.Range("J22").Value = 485
$ws.Range("K22").Value = 308
$ws.Range("L22").Value = 485
$ws.Range("M22").Value = 42
$ws.Range("N22").Value = -1185
$ws.Range("H24").Value = 4181.6
$ws.Range("I24").Value = 225.5
$ws.Range("J24").Value = 20006
$ws.Range("K24").Value = 225.5
$ws.Range("L24").Value = 20006
$ws.Range("M24").Value = -55.5
$ws.Range("N24").Value = -20346
$ws.Range("H31").Value = 2252
$ws.Range("I31").Value = 2142.4
$ws.Range("K31").Value = 2142.4
$ws.Range("M31").Value = -1847.4
$ws.Range("H34").Value = 2252
$ws.Range("I34").Value = 2142.4
$ws.Range("K34").Value = 2142.4
$ws.Range("M34").Value = -1940.4
$ws.Range("H62").Value = 2466.6667
$ws.Range("J62").Value = 2500
$ws.Range("L62").Value = 2500
$ws.Range("N62").Value = -3748
$ws.Range("H65").Value = 2466.6667
$ws.Range("J65").Value = 2500
$ws.Range("L65").Value = 12500
$ws.Range("N65").Value = -18740
$ws.Range("H86").Value = 4002.8
$ws.Range("I86").Value = 4003.5
$ws.Range("K86").Value = 4003.5
$ws.Range("M86").Value = -2880.5
$ws.Range("H89").Value = 4002.8
$ws.Range("I89").Value = 4003.5
$ws.Range("K89").Value = 20017.5
$ws.Range("M89").Value = -14401.5
$ws.Range("H107").Value = 1052.6666
$ws.Range("I107").Value = 532.6667
$ws.Range("J107").Value = 1572.6666
$ws.Range("K107").Value = 532.6667
$ws.Range("L107").Value = 1572.6666
$ws.Range("M107").Value = 1387.3333
$ws.Range("N107").Value = -5412.6666
$ws.Range("H125").Value = 99999
$ws.Range("J125").Value = 99999
$ws.Range("L125").Value = 99999
$ws.Range("N125").Value = -104919
$ws.Range("H141").Value = 467226.88
$ws.Range("J141").Value = 467226.88
$ws.Range("L141").Value = 467226.88
$ws.Range("N141").Value = -477586.88

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 44.458332
$ws.Range("I2").Value = 43.77778
$ws.Range("J2").Value = 44.866665
$ws.Range("K2").Value = 262.66668
$ws.Range("L2").Value = 269.19999
$ws.Range("M2").Value = -149.66668
$ws.Range("N2").Value = -495.19999
$ws.Range("H33").Value = 166.33333
$ws.Range("I33").Value = 150
$ws.Range("J33").Value = 199
$ws.Range("K33").Value = 900
$ws.Range("L33").Value = 1194
$ws.Range("M33").Value = -617
$ws.Range("N33").Value = -1760
$ws.Range("H52").Value = 20000
$ws.Range("J52").Value = 20000
$ws.Range("L52").Value = 60000
$ws.Range("N52").Value = -60532

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6864.7144
$ws.Range("I70").Value = 6864.7144
$ws.Range("K70").Value = 6864.7144
$ws.Range("M70").Value = -6594.7144
$ws.Range("H73").Value = 6864.7144
$ws.Range("I73").Value = 6864.7144
$ws.Range("K73").Value = 6864.7144
$ws.Range("M73").Value = -5928.7144
$ws.Range("H80").Value = 1800
$ws.Range("I80").Value = 1800
$ws.Range("K80").Value = 1800
$ws.Range("M80").Value = -802
$ws.Range("H83").Value = 1800
$ws.Range("I83").Value = 1800
$ws.Range("K83").Value = 9000
$ws.Range("M83").Value = -4008

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H55").Value = 1781.9333
$ws.Range("I55").Value = 1671.6
$ws.Range("J55").Value = 1837.1
$ws.Range("K55").Value = 1671.6
$ws.Range("L55").Value = 1837.1
$ws.Range("M55").Value = -1498.6
$ws.Range("N55").Value = -2183.1
$ws.Range("H68").Value = 1894.3334
$ws.Range("I68").Value = 1098
$ws.Range("K68").Value = 1098
$ws.Range("M68").Value = -349
$ws.Range("H71").Value = 1894.3334
$ws.Range("I71").Value = 1098
$ws.Range("K71").Value = 5490
$ws.Range("M71").Value = -1746
